$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Creación del servidor de Discord (medio de comunicación principal)."
#    Remove the spell-check proofErr split around "Discord" by collapsing the
#    sentence's runs back into one (text itself is unchanged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Creación del servidor de Discord (medio de comunicación principal).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Creación del servidor de Discord (medio de comunicación principal).", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "... de tomar contacto con cada uno de los roles y subequipos que tiene
#     cada miembro asignado. " - remove proofErr split around "subequipos".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "de tomar contacto con cada uno de los roles y subequipos que tiene cada miembro asignado. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "de tomar contacto con cada uno de los roles y subequipos que tiene cada miembro asignado. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Inauguración del tablero de actividades en GitHub" -> add trailing period
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Inauguración del tablero de actividades en GitHub",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Inauguración del tablero de actividades en GitHub.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "[Falta de asistencia] Daniel Yanel Gorrón y Rubén López"
#    -> "[Falta de asistencia injustificada] Daniel Yanel Gorrón y Rubén López."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "[Falta de asistencia] Daniel Yanel Gorrón y Rubén López",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Falta de asistencia injustificada] Daniel Yanel Gorrón y Rubén López.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Append the new "Sexta Reunión" section as new paragraphs at the end of
#    the document. Build every paragraph first with plain/default formatting,
#    and only apply the bold/underline/size-32 heading look afterwards so the
#    formatting does not leak forward into paragraphs created later.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$pBlank = $d.Paragraphs.Last
# Blank centered paragraph - nothing else to do.

$pBlank.Range.InsertParagraphAfter()
$pHeading = $d.Paragraphs.Last
$pHeading.Range.Text = "Sexta Reunión (12/04/2022) "

$pHeading.Range.InsertParagraphAfter()
$pFase = $d.Paragraphs.Last
$pFase.Range.Text = "Inicio de la fase 2 de la planificación de desarrollo del proyecto."

$pFase.Range.InsertParagraphAfter()
$pFalta1 = $d.Paragraphs.Last
$pFalta1.Range.Text = "[Falta de asistencia injustificada] Iubal Nicolás Camjalli Spiegel, Mario González Montalvo, Daniel Yanel Gorrón, Rubén López, Fátima De la Morena y Carlos Jiménez Crespo."

$pFalta1.Range.InsertParagraphAfter()
$pFalta2 = $d.Paragraphs.Last
$pFalta2.Range.Text = "[Falta de realización de actividades injustificada]  Daniel Yanel Gorrón y Rubén López."

# Now apply the bold / underlined / size-32 heading formatting to the
# "Sexta Reunión" paragraph only (done last so later paragraphs already
# exist and keep their plain formatting).
$pHeading.Range.Font.Bold = $true
$pHeading.Range.Font.BoldBi = $true
$pHeading.Range.Font.Underline = 1
$pHeading.Range.Font.Size = 16
$pHeading.Range.Font.SizeBi = 16

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
